$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) rows 2-6 from 2023-10-05 (45204) to 2023-10-08 (45207)
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = 45207
}
